# Atualiza dados da liga eliminacao
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New time_id values (column A), rows 13..25, shifted up by one
# row 13 ("BORGES ITAQUI F.C.") leaves the table, every other team's
# id moves up one row, and a brand new team/id pair is appended at
# the bottom (row 25).
$ids = @(5823700, 13707047, 13913874, 13951133, 14124559, 18223508, 18344271, 18642587, 18661583, 19033717, 20696550, 24468241, 24856400)

# New team names (column B) for rows 13..25, same shift-up pattern,
# with "Grêmio imortal 36" newly introduced as the last entry.
$names = @("S.E.R. GRILLO", "Super Vasco f.c", "Bandoleros FCS", "JUV. KP", "Paulo Virgili FC", "Rolo Compressor  ZN", "FÚRIA LEON", "Fedato Futebol Clube", "pura bucha/internacional", "Mau Humor F.C.", "Dom Camillo68", "Grêmio imortal 37", "Grêmio imortal 36")

for ($i = 0; $i -lt $ids.Length; $i++) {
    $row = 13 + $i
    $ws.Cells.Item($row, 1).Value = $ids[$i]
    $ws.Cells.Item($row, 2).Value = $names[$i]
}
